$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test-Cases")

# Row 7: change Approved/Rejected status to "Rejected" and add a reason
$ws.Range("I7").Value = "Rejected"
$ws.Range("J7").Value = "trest"

# Row 10: change Approved/Rejected status to "Rejected" and add a reason
$ws.Range("I10").Value = "Rejected"
$ws.Range("J10").Value = "tesrt"

# Update selection to reflect the last edited cell
$ws.Range("I10").Select()
